$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 13 (Squilla mantis, station's 1-RAP gear) weight/number values
$ws.Range("G13").Value = 0.038
$ws.Range("H13").Value = 2

# Remove the duplicate "Squilla mantis" row (row 25) for the 2-RAP gear block;
# this shifts rows 26:43 up to become rows 25:42
$ws.Rows.Item(25).Delete()

# After the shift, update Numb for the three rows whose count became -1
$ws.Range("H32").Value = -1
$ws.Range("H39").Value = -1
$ws.Range("H42").Value = -1

# Recompute RF for the shifted block (rows 25:42) to the new factor
$ws.Range("I25:I42").Value = 4.807704918032787
